$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 33 with new time entry
$ws.Cells.Item(33, 1).Value = 44056
$ws.Cells.Item(33, 2).Value = 4
$ws.Cells.Item(33, 3).Value = "Kuvien lataaminen serveriltä sekä productio version testaus"

# Set row height for row 33 to match the other wrapped-text rows (30pt)
$ws.Rows.Item(33).RowHeight = 30

# Update the active selection to reflect the edited cell
$ws.Range("C33").Select()
